$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A8: "Volume 30   Number  9" -> "Volume 30   Number  10" ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 1).Text = "10"

# --- C9: "Report Covering the Week  2/27/2023  Through  3/5/2023"
#         -> "Report Covering the Week  3/6/2023  Through  3/12/2023" ---
# Update the 2nd (later) date first so the 1st date's character offsets stay valid.
$c9 = $ws.Range("C9")
$c9.Characters(47, 8).Text = "3/12/2023"
$c9.Characters(27, 9).Text = "3/6/2023"

# --- Weekly crime-stat numbers (rows 14-30) ---
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = -33.333333333333
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = 12
$ws.Range("H14").Value = -41.666666666666
$ws.Range("I14").Value = 20
$ws.Range("J14").Value = 25
$ws.Range("K14").Value = -20
$ws.Range("L14").Value = -4.761904761904
$ws.Range("M14").Value = 11.111111111111
$ws.Range("N14").Value = -77.52808988764
$ws.Range("C15").Value = 12
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = 50
$ws.Range("G15").Value = 27
$ws.Range("H15").Value = 7.407407407407
$ws.Range("I15").Value = 76
$ws.Range("J15").Value = 81
$ws.Range("K15").Value = -6.172839506172
$ws.Range("L15").Value = 16.923076923076
$ws.Range("M15").Value = 55.102040816326
$ws.Range("N15").Value = -30.90909090909
$ws.Range("C16").Value = 87
$ws.Range("D16").Value = 83
$ws.Range("E16").Value = 4.819277108433
$ws.Range("F16").Value = 302
$ws.Range("G16").Value = 348
$ws.Range("H16").Value = -13.218390804597
$ws.Range("I16").Value = 820
$ws.Range("J16").Value = 836
$ws.Range("K16").Value = -1.913875598086
$ws.Range("L16").Value = 33.986928104575
$ws.Range("M16").Value = 10.81081081081
$ws.Range("N16").Value = -74.28661022264
$ws.Range("C17").Value = 128
$ws.Range("D17").Value = 126
$ws.Range("E17").Value = 1.587301587301
$ws.Range("F17").Value = 501
$ws.Range("G17").Value = 459
$ws.Range("H17").Value = 9.150326797385
$ws.Range("I17").Value = 1285
$ws.Range("J17").Value = 1171
$ws.Range("K17").Value = 9.735269000853
$ws.Range("L17").Value = 25.243664717348
$ws.Range("M17").Value = 62.247474747474
$ws.Range("N17").Value = -11.805078929306
$ws.Range("C18").Value = 65
$ws.Range("D18").Value = 74
$ws.Range("E18").Value = -12.162162162162
$ws.Range("F18").Value = 224
$ws.Range("G18").Value = 254
$ws.Range("H18").Value = -11.811023622047
$ws.Range("I18").Value = 578
$ws.Range("J18").Value = 579
$ws.Range("K18").Value = -0.172711571675
$ws.Range("L18").Value = 44.5
$ws.Range("M18").Value = -4.145936981757
$ws.Range("N18").Value = -83.732057416267
$ws.Range("C19").Value = 125
$ws.Range("D19").Value = 167
$ws.Range("E19").Value = -25.149700598802
$ws.Range("F19").Value = 524
$ws.Range("G19").Value = 630
$ws.Range("H19").Value = -16.825396825396
$ws.Range("I19").Value = 1344
$ws.Range("J19").Value = 1501
$ws.Range("K19").Value = -10.459693537641
$ws.Range("L19").Value = 29.230769230769
$ws.Range("M19").Value = 79.919678714859
$ws.Range("N19").Value = 4.510108864696
$ws.Range("C20").Value = 94
$ws.Range("D20").Value = 82
$ws.Range("E20").Value = 14.634146341463
$ws.Range("F20").Value = 369
$ws.Range("G20").Value = 333
$ws.Range("H20").Value = 10.81081081081
$ws.Range("I20").Value = 1019
$ws.Range("J20").Value = 892
$ws.Range("K20").Value = 14.237668161435
$ws.Range("L20").Value = 163.989637305699
$ws.Range("M20").Value = 179.945054945055
$ws.Range("N20").Value = -66.067266067266
$ws.Range("C21").Value = 513
$ws.Range("D21").Value = 543
$ws.Range("E21").Value = -5.524861878453
$ws.Range("F21").Value = 1956
$ws.Range("G21").Value = 2063
$ws.Range("H21").Value = -5.186621425109
$ws.Range("I21").Value = 5142
$ws.Range("J21").Value = 5085
$ws.Range("K21").Value = 1.120943952802
$ws.Range("L21").Value = 44.845070422535
$ws.Range("M21").Value = 55.206761243585
$ws.Range("N21").Value = -59.470323953653
$ws.Range("C22").Value = 6
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 29
$ws.Range("H22").Value = -13.793103448275
$ws.Range("I22").Value = 48
$ws.Range("J22").Value = 71
$ws.Range("K22").Value = -32.394366197183
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = 9.090909090909
$ws.Range("C23").Value = 21
$ws.Range("D23").Value = 29
$ws.Range("E23").Value = -27.586206896551
$ws.Range("F23").Value = 113
$ws.Range("G23").Value = 107
$ws.Range("H23").Value = 5.607476635514
$ws.Range("I23").Value = 315
$ws.Range("J23").Value = 276
$ws.Range("K23").Value = 14.130434782608
$ws.Range("L23").Value = 52.173913043478
$ws.Range("M23").Value = 75.977653631284
$ws.Range("C24").Value = 291
$ws.Range("D24").Value = 335
$ws.Range("E24").Value = -13.134328358209
$ws.Range("F24").Value = 1330
$ws.Range("G24").Value = 1388
$ws.Range("H24").Value = -4.178674351585
$ws.Range("I24").Value = 3189
$ws.Range("J24").Value = 3176
$ws.Range("K24").Value = 0.409319899244
$ws.Range("L24").Value = 33.766778523489
$ws.Range("M24").Value = 46.823204419889
$ws.Range("C25").Value = 170
$ws.Range("D25").Value = 167
$ws.Range("E25").Value = 1.796407185628
$ws.Range("F25").Value = 761
$ws.Range("G25").Value = 703
$ws.Range("H25").Value = 8.250355618776
$ws.Range("I25").Value = 1859
$ws.Range("J25").Value = 1735
$ws.Range("K25").Value = 7.1469740634
$ws.Range("L25").Value = 30.823363828289
$ws.Range("M25").Value = 1.032608695652
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -18.75
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 55
$ws.Range("H26").Value = -16.363636363636
$ws.Range("I26").Value = 123
$ws.Range("J26").Value = 142
$ws.Range("K26").Value = -13.38028169014
$ws.Range("L26").Value = 20.588235294117
$ws.Range("D27").Value = 18
$ws.Range("E27").Value = 27.777777777777
$ws.Range("F27").Value = 82
$ws.Range("G27").Value = 71
$ws.Range("H27").Value = 15.492957746478
$ws.Range("I27").Value = 208
$ws.Range("J27").Value = 160
$ws.Range("K27").Value = 30
$ws.Range("L27").Value = 30.817610062893
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 8
$ws.Range("E28").Value = -37.5
$ws.Range("F28").Value = 18
$ws.Range("G28").Value = 22
$ws.Range("H28").Value = -18.181818181818
$ws.Range("I28").Value = 56
$ws.Range("J28").Value = 80
$ws.Range("K28").Value = -30
$ws.Range("L28").Value = -8.196721311475
$ws.Range("M28").Value = -21.12676056338
$ws.Range("N28").Value = -72.549019607843
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 8
$ws.Range("E29").Value = -37.5
$ws.Range("F29").Value = 17
$ws.Range("G29").Value = 22
$ws.Range("H29").Value = -22.727272727272
$ws.Range("I29").Value = 46
$ws.Range("J29").Value = 75
$ws.Range("K29").Value = -38.666666666666
$ws.Range("L29").Value = -17.857142857142
$ws.Range("M29").Value = -23.333333333333
$ws.Range("N29").Value = -74.725274725274
$ws.Range("G30").Value = 6
